$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 128
$ws.Cells.Item(128, 3).Value = 'Infernal Thorn'
$ws.Cells.Item(128, 4).Value = 'dagger'
$ws.Cells.Item(128, 6).Value = 'Hell Forged'
$ws.Cells.Item(128, 7).Value = 'Forged in the fires of the abyss, this dagger sears the soul of its victims.'
$ws.Cells.Item(128, 9).Value = 2000
$ws.Cells.Item(128, 12).Value = 150000000000
$ws.Cells.Item(128, 13).Value = 10000
$ws.Cells.Item(128, 14).Value = 1000
$ws.Cells.Item(128, 15).Value = 0
$ws.Cells.Item(128, 17).Value = 0.1
$ws.Cells.Item(128, 22).Value = 0.45
$ws.Cells.Item(128, 25).Value = 0
$ws.Cells.Item(128, 26).Value = 0
$ws.Cells.Item(128, 29).Value = 1
$ws.Cells.Item(128, 35).Value = 0
$ws.Cells.Item(128, 36).Value = 0
$ws.Cells.Item(128, 37).Value = 0
$ws.Cells.Item(128, 38).Value = 0
$ws.Cells.Item(128, 39).Value = 0
$ws.Cells.Item(128, 45).Value = 0
$ws.Cells.Item(128, 48).Value = 0
$ws.Cells.Item(128, 50).Value = 0
$ws.Cells.Item(128, 51).Value = 0
$ws.Cells.Item(128, 53).Value = 0
$ws.Cells.Item(128, 54).Value = 0
$ws.Cells.Item(128, 55).Value = 0
$ws.Cells.Item(128, 56).Value = 0
$ws.Cells.Item(128, 57).Value = 0
$ws.Cells.Item(128, 58).Value = 0
$ws.Cells.Item(128, 59).Value = 0
$ws.Cells.Item(128, 65).Value = 20
$ws.Cells.Item(128, 66).Value = 0
$ws.Cells.Item(128, 67).Value = 0
$ws.Cells.Item(128, 68).Value = 0
$ws.Cells.Item(128, 69).Value = 0

# Row 129
$ws.Cells.Item(129, 3).Value = 'Blackwater Shiv'
$ws.Cells.Item(129, 4).Value = 'dagger'
$ws.Cells.Item(129, 6).Value = 'Pirate Lord Leather'
$ws.Cells.Item(129, 7).Value = 'A deadly blade carried by the infamous Pirate Lords, swift and merciless.'
$ws.Cells.Item(129, 9).Value = 5000
$ws.Cells.Item(129, 12).Value = 1000000000
$ws.Cells.Item(129, 15).Value = 0
$ws.Cells.Item(129, 17).Value = 0.3
$ws.Cells.Item(129, 22).Value = 0.7
$ws.Cells.Item(129, 25).Value = 0
$ws.Cells.Item(129, 26).Value = 0
$ws.Cells.Item(129, 29).Value = 1
$ws.Cells.Item(129, 35).Value = 0
$ws.Cells.Item(129, 36).Value = 0
$ws.Cells.Item(129, 37).Value = 0
$ws.Cells.Item(129, 38).Value = 0
$ws.Cells.Item(129, 39).Value = 0
$ws.Cells.Item(129, 45).Value = 0
$ws.Cells.Item(129, 48).Value = 0
$ws.Cells.Item(129, 50).Value = 0
$ws.Cells.Item(129, 51).Value = 0
$ws.Cells.Item(129, 53).Value = 0
$ws.Cells.Item(129, 54).Value = 0
$ws.Cells.Item(129, 55).Value = 0
$ws.Cells.Item(129, 56).Value = 0
$ws.Cells.Item(129, 57).Value = 0
$ws.Cells.Item(129, 58).Value = 0
$ws.Cells.Item(129, 59).Value = 0
$ws.Cells.Item(129, 65).Value = 20
$ws.Cells.Item(129, 66).Value = 0
$ws.Cells.Item(129, 67).Value = 0
$ws.Cells.Item(129, 68).Value = 0
$ws.Cells.Item(129, 69).Value = 0

# Row 130
$ws.Cells.Item(130, 3).Value = 'Soulpiercer'
$ws.Cells.Item(130, 4).Value = 'dagger'
$ws.Cells.Item(130, 6).Value = 'Purgatory Chains'
$ws.Cells.Item(130, 7).Value = 'Enchanted with spectral chains, this dagger binds the souls of the fallen.'
$ws.Cells.Item(130, 9).Value = 3000
$ws.Cells.Item(130, 12).Value = 250000000000
$ws.Cells.Item(130, 13).Value = 25000
$ws.Cells.Item(130, 14).Value = 2500
$ws.Cells.Item(130, 15).Value = 1000
$ws.Cells.Item(130, 17).Value = 0.15
$ws.Cells.Item(130, 22).Value = 0.55
$ws.Cells.Item(130, 25).Value = 0
$ws.Cells.Item(130, 26).Value = 0
$ws.Cells.Item(130, 29).Value = 1
$ws.Cells.Item(130, 35).Value = 0
$ws.Cells.Item(130, 36).Value = 0
$ws.Cells.Item(130, 37).Value = 0
$ws.Cells.Item(130, 38).Value = 0
$ws.Cells.Item(130, 39).Value = 0
$ws.Cells.Item(130, 45).Value = 0
$ws.Cells.Item(130, 48).Value = 0
$ws.Cells.Item(130, 50).Value = 0
$ws.Cells.Item(130, 51).Value = 0
$ws.Cells.Item(130, 53).Value = 0
$ws.Cells.Item(130, 54).Value = 0
$ws.Cells.Item(130, 55).Value = 0
$ws.Cells.Item(130, 56).Value = 0
$ws.Cells.Item(130, 57).Value = 0
$ws.Cells.Item(130, 58).Value = 0
$ws.Cells.Item(130, 59).Value = 0
$ws.Cells.Item(130, 65).Value = 20
$ws.Cells.Item(130, 66).Value = 0
$ws.Cells.Item(130, 67).Value = 0
$ws.Cells.Item(130, 68).Value = 0
$ws.Cells.Item(130, 69).Value = 0

# Row 131
$ws.Cells.Item(131, 3).Value = 'Frostbite Fang'
$ws.Cells.Item(131, 4).Value = 'dagger'
$ws.Cells.Item(131, 6).Value = 'Corrupted Ice'
$ws.Cells.Item(131, 7).Value = 'A shard of frozen corruption, it freezes the hearts of those it cuts.'
$ws.Cells.Item(131, 9).Value = 5000
$ws.Cells.Item(131, 12).Value = 1000000000
$ws.Cells.Item(131, 15).Value = 0
$ws.Cells.Item(131, 17).Value = 0.3
$ws.Cells.Item(131, 22).Value = 0.75
$ws.Cells.Item(131, 25).Value = 0
$ws.Cells.Item(131, 26).Value = 0
$ws.Cells.Item(131, 29).Value = 1
$ws.Cells.Item(131, 35).Value = 0
$ws.Cells.Item(131, 36).Value = 0
$ws.Cells.Item(131, 37).Value = 0
$ws.Cells.Item(131, 38).Value = 0
$ws.Cells.Item(131, 39).Value = 0
$ws.Cells.Item(131, 45).Value = 0
$ws.Cells.Item(131, 48).Value = 0
$ws.Cells.Item(131, 50).Value = 0
$ws.Cells.Item(131, 51).Value = 0
$ws.Cells.Item(131, 53).Value = 0
$ws.Cells.Item(131, 54).Value = 0
$ws.Cells.Item(131, 55).Value = 0
$ws.Cells.Item(131, 56).Value = 0
$ws.Cells.Item(131, 57).Value = 0
$ws.Cells.Item(131, 58).Value = 0
$ws.Cells.Item(131, 59).Value = 0
$ws.Cells.Item(131, 65).Value = 20
$ws.Cells.Item(131, 66).Value = 0
$ws.Cells.Item(131, 67).Value = 0
$ws.Cells.Item(131, 68).Value = 0
$ws.Cells.Item(131, 69).Value = 0

# Row 132
$ws.Cells.Item(132, 3).Value = 'Stonefang Dagger'
$ws.Cells.Item(132, 4).Value = 'dagger'
$ws.Cells.Item(132, 6).Value = 'Twisted Earth'
$ws.Cells.Item(132, 7).Value = 'Formed from the depths of the shifting earth, it trembles with untamed power.'
$ws.Cells.Item(132, 9).Value = 5000
$ws.Cells.Item(132, 12).Value = 500000000000
$ws.Cells.Item(132, 13).Value = 50000
$ws.Cells.Item(132, 14).Value = 5000
$ws.Cells.Item(132, 15).Value = 10000
$ws.Cells.Item(132, 17).Value = 0.3
$ws.Cells.Item(132, 22).Value = 0.75
$ws.Cells.Item(132, 25).Value = 0
$ws.Cells.Item(132, 26).Value = 0
$ws.Cells.Item(132, 29).Value = 1
$ws.Cells.Item(132, 35).Value = 0
$ws.Cells.Item(132, 36).Value = 0
$ws.Cells.Item(132, 37).Value = 0
$ws.Cells.Item(132, 38).Value = 0
$ws.Cells.Item(132, 39).Value = 0
$ws.Cells.Item(132, 45).Value = 0
$ws.Cells.Item(132, 48).Value = 0
$ws.Cells.Item(132, 50).Value = 0
$ws.Cells.Item(132, 51).Value = 0
$ws.Cells.Item(132, 53).Value = 0
$ws.Cells.Item(132, 54).Value = 0
$ws.Cells.Item(132, 55).Value = 0
$ws.Cells.Item(132, 56).Value = 0
$ws.Cells.Item(132, 57).Value = 0
$ws.Cells.Item(132, 58).Value = 0
$ws.Cells.Item(132, 59).Value = 0
$ws.Cells.Item(132, 65).Value = 20
$ws.Cells.Item(132, 66).Value = 0
$ws.Cells.Item(132, 67).Value = 0
$ws.Cells.Item(132, 68).Value = 0
$ws.Cells.Item(132, 69).Value = 0

# Row 133
$ws.Cells.Item(133, 3).Value = 'Phantom''s Veil'
$ws.Cells.Item(133, 4).Value = 'dagger'
$ws.Cells.Item(133, 6).Value = 'Delusional Silver'
$ws.Cells.Item(133, 7).Value = 'A blade of shimmering deception, warping reality with each strike.'
$ws.Cells.Item(133, 9).Value = 8000
$ws.Cells.Item(133, 12).Value = 1000000000
$ws.Cells.Item(133, 15).Value = 0
$ws.Cells.Item(133, 17).Value = 0.4
$ws.Cells.Item(133, 22).Value = 0.9
$ws.Cells.Item(133, 25).Value = 0
$ws.Cells.Item(133, 26).Value = 0
$ws.Cells.Item(133, 29).Value = 1
$ws.Cells.Item(133, 35).Value = 0
$ws.Cells.Item(133, 36).Value = 0
$ws.Cells.Item(133, 37).Value = 0
$ws.Cells.Item(133, 38).Value = 0
$ws.Cells.Item(133, 39).Value = 0
$ws.Cells.Item(133, 45).Value = 0
$ws.Cells.Item(133, 48).Value = 0
$ws.Cells.Item(133, 50).Value = 0
$ws.Cells.Item(133, 51).Value = 0
$ws.Cells.Item(133, 53).Value = 0
$ws.Cells.Item(133, 54).Value = 0
$ws.Cells.Item(133, 55).Value = 0
$ws.Cells.Item(133, 56).Value = 0
$ws.Cells.Item(133, 57).Value = 0
$ws.Cells.Item(133, 58).Value = 0
$ws.Cells.Item(133, 59).Value = 0
$ws.Cells.Item(133, 65).Value = 20
$ws.Cells.Item(133, 66).Value = 0
$ws.Cells.Item(133, 67).Value = 0
$ws.Cells.Item(133, 68).Value = 0
$ws.Cells.Item(133, 69).Value = 0

# Row 134
$ws.Cells.Item(134, 3).Value = 'Oathsever Blade'
$ws.Cells.Item(134, 4).Value = 'dagger'
$ws.Cells.Item(134, 6).Value = 'Faithless Plate'
$ws.Cells.Item(134, 7).Value = 'Once wielded by fallen knights, this dagger carries the weight of broken oaths.'
$ws.Cells.Item(134, 9).Value = 8000
$ws.Cells.Item(134, 12).Value = 1000000000
$ws.Cells.Item(134, 15).Value = 0
$ws.Cells.Item(134, 17).Value = 0.4
$ws.Cells.Item(134, 22).Value = 0.9
$ws.Cells.Item(134, 25).Value = 0
$ws.Cells.Item(134, 26).Value = 0
$ws.Cells.Item(134, 29).Value = 1
$ws.Cells.Item(134, 35).Value = 0
$ws.Cells.Item(134, 36).Value = 0
$ws.Cells.Item(134, 37).Value = 0
$ws.Cells.Item(134, 38).Value = 0
$ws.Cells.Item(134, 39).Value = 0
$ws.Cells.Item(134, 45).Value = 0
$ws.Cells.Item(134, 48).Value = 0
$ws.Cells.Item(134, 50).Value = 0
$ws.Cells.Item(134, 51).Value = 0
$ws.Cells.Item(134, 53).Value = 0
$ws.Cells.Item(134, 54).Value = 0
$ws.Cells.Item(134, 55).Value = 0
$ws.Cells.Item(134, 56).Value = 0
$ws.Cells.Item(134, 57).Value = 0
$ws.Cells.Item(134, 58).Value = 0
$ws.Cells.Item(134, 59).Value = 0
$ws.Cells.Item(134, 65).Value = 20
$ws.Cells.Item(134, 66).Value = 0
$ws.Cells.Item(134, 67).Value = 0
$ws.Cells.Item(134, 68).Value = 0
$ws.Cells.Item(134, 69).Value = 0

# Fix row 114 height (was 23.85, should now be default 13.8)
$ws.Rows.Item(114).RowHeight = 13.8

# Update sheet view: scroll position and active cell selection
$win = $excel.ActiveWindow
$win.ScrollRow = 113
$win.ScrollColumn = 1
$ws.Range("A134").Select()
